$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This document has two identical footnote-style paragraphs:
#   "*The National Percentile Rank is not calculated for Autumn tests."
# which need to become:
#   "*The National Percentile Rank is not calculated for the Autumn test."
# split across three runs (the way Word splits runs when the sentence is
# edited in two separate spots), and the document's lone "_GoBack" bookmark
# (an artifact of the last cursor position when the file was last saved)
# needs to move from its old location to sit between run #2 and run #3 of
# the *second* occurrence of that sentence.
# ---------------------------------------------------------------------------

# Drop the old/stale "_GoBack" bookmark, wherever it currently sits - Word
# re-anchors this automatically to wherever text was most recently edited.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Collapse a Range to a single point and add+delete a throw-away bookmark
# there. Word's bookmark engine splits the run under the bookmark's
# start/end, so this forces a run break at that character offset without
# altering the surrounding text or its formatting.
function Split-RunAt($point) {
    $null = $d.Bookmarks.Add("__tmp_split__", $d.Range($point, $point))
    $d.Bookmarks("__tmp_split__").Delete()
}

# Rewrites one occurrence of "...calculated for Autumn tests." into three
# runs: "...calculated for" | " the Autumn test" | "." -- optionally
# stitching the real "_GoBack" bookmark between run #2 and run #3.
function Update-Sentence($searchFrom, [bool]$placeBookmark) {
    $found = $d.Range($searchFrom, $d.Content.End)
    $null = $found.Find.Execute("for Autumn tests.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $forEnd = $found.Start + 3      # right after "for"
    $sentenceEnd = $found.End       # right after the final "."
    $sPos = $sentenceEnd - 2        # the "s" in "tests", right before "."

    # delete the trailing "s" of "tests" first (rightmost edit), then
    # insert "the " right after "for " (leftmost edit) so the earlier
    # offset is untouched by the first deletion.
    $d.Range($sPos, $sPos + 1).Text = ""
    $d.Range($forEnd + 1, $forEnd + 1).Text = "the "

    $newSentenceEnd = $sentenceEnd + 4 - 1   # +4 for "the ", -1 for removed "s"
    $periodPos = $newSentenceEnd - 1

    if ($placeBookmark) {
        Split-RunAt $periodPos
        $null = $d.Bookmarks.Add("_GoBack", $d.Range($periodPos, $periodPos))
    } else {
        Split-RunAt $periodPos
    }
    Split-RunAt $forEnd

    return $newSentenceEnd
}

$end1 = Update-Sentence 0 $false
$null = Update-Sentence $end1 $true

Write-Output "done"
